$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.043.06'
$ws.Range('E2').Value = '  +1.50%  '
$ws.Range('D3').Value = '1.712.73'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").ClearFormats()
$ws.Range('E4').Value = '  +0.81%  '
$ws.Range("D5").Value = "'318.37"
$ws.Range("D5").ClearFormats()
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range("D6").Value = "'1.008"
$ws.Range("D6").ClearFormats()
$ws.Range('E6').Value = '  +0.78%  '
$ws.Range("D7").Value = "'0.3963"
$ws.Range("D7").ClearFormats()
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range("D8").Value = "'0.4114"
$ws.Range("D8").ClearFormats()
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range("D9").Value = "'1.526"
$ws.Range("D9").ClearFormats()
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('E10').Value = '  +0.80%  '
$ws.Range("D11").Value = "'52.44"
$ws.Range("D11").ClearFormats()
$ws.Range('E11').Value = '  +2.17%  '
$ws.Range("D12").Value = "'0.08897"
$ws.Range("D12").ClearFormats()
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range("D13").Value = "'7.687"
$ws.Range("D13").ClearFormats()
$ws.Range('E13').Value = '  +6.74%  '
$ws.Range("D14").Value = "'24.92"
$ws.Range("D14").ClearFormats()
$ws.Range('E14').Value = '  +4.47%  '
$ws.Range("D15").Value = "'0.00001390"
$ws.Range("D15").ClearFormats()
$ws.Range('E15').Value = '  +4.35%  '
$ws.Range("D16").Value = "'8.124"
$ws.Range("D16").ClearFormats()
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('D17').Value = '1.712.27'
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range("D18").Value = "'100.58"
$ws.Range("D18").ClearFormats()
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range("D19").Value = "'0.07126"
$ws.Range("D19").ClearFormats()
$ws.Range('E19').Value = '  +1.56%  '
$ws.Range("D20").Value = "'20.16"
$ws.Range("D20").ClearFormats()
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range("D21").Value = "'7.456"
$ws.Range("D21").ClearFormats()
$ws.Range('E21').Value = '  +5.98%  '
$ws.Range('E22').Value = '  +0.84%  '
$ws.Range("D23").Value = "'14.46"
$ws.Range("D23").ClearFormats()
$ws.Range('E23').Value = '  +1.31%  '
$ws.Range('D24').Value = '25.031.89'
$ws.Range('E24').Value = '  +1.43%  '
$ws.Range("D25").Value = "'3.092"
$ws.Range("D25").ClearFormats()
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range("D26").Value = "'2.353"
$ws.Range("D26").ClearFormats()
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range("D27").Value = "'23.01"
$ws.Range("D27").ClearFormats()
$ws.Range('E27').Value = '  +1.37%  '
$ws.Range("D28").Value = "'165.29"
$ws.Range("D28").ClearFormats()
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range("D29").Value = "'8.740"
$ws.Range("D29").ClearFormats()
$ws.Range('E29').Value = '  +17.20%  '
$ws.Range("D30").Value = "'138.83"
$ws.Range("D30").ClearFormats()
$ws.Range('E30').Value = '  +1.19%  '
$ws.Range("D31").Value = "'5.204"
$ws.Range("D31").ClearFormats()
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range("D32").Value = "'7.762"
$ws.Range("D32").ClearFormats()
$ws.Range('E32').Value = '  +8.64%  '
$ws.Range('D33').Value = '1.900.14'
$ws.Range('E33').Value = '  +1.06%  '
$ws.Range("D34").Value = "'0.08979"
$ws.Range("D34").ClearFormats()
$ws.Range('E34').Value = '  +4.21%  '
$ws.Range("D35").Value = "'1.072"
$ws.Range("D35").ClearFormats()
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range("D36").Value = "'0.02988"
$ws.Range("D36").ClearFormats()
$ws.Range('E36').Value = '  +9.42%  '
$ws.Range('B37').Value = 'WEMIXTOKEN'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = "'1.982"
$ws.Range("D37").ClearFormats()
$ws.Range('E37').Value = '  +2.89%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = "'0.2773"
$ws.Range("D38").ClearFormats()
$ws.Range('E38').Value = '  +0.95%  '
$ws.Range("D39").Value = "'11.04"
$ws.Range("D39").ClearFormats()
$ws.Range('E39').Value = '  -3.79%  '
$ws.Range("D40").Value = "'14.61"
$ws.Range("D40").ClearFormats()
$ws.Range('E40').Value = '  +0.82%  '
$ws.Range("D41").Value = "'0.8143"
$ws.Range("D41").ClearFormats()
$ws.Range('E41').Value = '  +5.90%  '
$ws.Range("D42").Value = "'0.09235"
$ws.Range("D42").ClearFormats()
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range("D43").Value = "'1.483"
$ws.Range("D43").ClearFormats()
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range("D44").Value = "'16.61"
$ws.Range("D44").ClearFormats()
$ws.Range('E44').Value = '  +3.91%  '
$ws.Range("D45").Value = "'0.7365"
$ws.Range("D45").ClearFormats()
$ws.Range('E45').Value = '  +2.44%  '
$ws.Range("D46").Value = "'2.632"
$ws.Range("D46").ClearFormats()
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range("D47").Value = "'4.287"
$ws.Range("D47").ClearFormats()
$ws.Range('E47').Value = '  +1.37%  '
$ws.Range("D48").Value = "'1.008"
$ws.Range("D48").ClearFormats()
$ws.Range('E48').Value = '  +0.77%  '
$ws.Range("D49").Value = "'1.341"
$ws.Range("D49").ClearFormats()
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range("D50").Value = "'140.14"
$ws.Range("D50").ClearFormats()
$ws.Range('E50').Value = '  -0.31%  '
$ws.Range("D51").Value = "'92.68"
$ws.Range("D51").ClearFormats()
$ws.Range('E51').Value = '  +3.48%  '
